$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49 - this shifts rows 49..133 down to 50..134,
# matching how the rest of the sheet already looks (same category/product
# columns repeated), and mirrors the row-carrying-row shift visible in the diff.
$ws.Rows("49:49").Insert()

# Populate the newly inserted row 49 with the new daily price record.
$ws.Range("A49").Value = 9
$ws.Range("B49").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44533
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100101
$ws.Range("H49").Value = "Berries"
$ws.Range("I49").Value = 100101001
$ws.Range("J49").Value = "Arándano (blue)"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 470
$ws.Range("N49").Value = 4000
$ws.Range("O49").Value = 5000
$ws.Range("P49").Value = 4468
$ws.Range("Q49").Value = "$/bandeja 2 kilos"
$ws.Range("R49").Value = "Provincia de Curicó"
$ws.Range("S49").Value = 2234
$ws.Range("T49").Value = 2
